$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C values based on which group (column B) each row belongs to.
# Group B=2 -> rows 2-92:   1.2E-2  => 2.3E-2
$ws.Range("C2:C92").Value = 0.023

# Group B=3 -> rows 93-254: 3.3000000000000002E-2 => 2.8000000000000001E-2
$ws.Range("C93:C254").Value = 0.028

# Group B=5 -> rows 255-268: 3.3000000000000002E-2 => 0.11
$ws.Range("C255:C268").Value = 0.11

# Group B=6 -> rows 269-314: 4.4999999999999998E-2 => 5.8999999999999997E-2
$ws.Range("C269:C314").Value = 0.059
